$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.252.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.427.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.00%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'555.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.12%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.80%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'Cardano"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.355"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'TRON"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.22%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'25.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.85%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.860.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.02%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'60.153.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.15%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.87%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.385.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.66%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'329.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.57%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'66.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.83%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.179"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.67%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.84%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0780"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.22%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'170.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.407"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.53%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'18.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").Value = "'4.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.84%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'332.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'145.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.06%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'InjectiveProtocol"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'20.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Stellar"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0969"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.578"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.91%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0225"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.07%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.04%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.28%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied changes"